# Documenting espn test cases:
#   - rename "Sheet3" -> "nbaCityNameTeamName"
#   - add a new "nbaPage" sheet (bold "Assertions" header + "NBA" note)
#   - on "nbaTeams", stamp an extra (empty, text-formatted) row under the
#     existing data and leave the selection sitting there

$wb = $excel.ActiveWorkbook

# --- rename the old "Sheet3" tab ---------------------------------------
$wsCity = $wb.Worksheets.Item("Sheet3")
$wsCity.Name = "nbaCityNameTeamName"

# --- append a brand-new "nbaPage" sheet after the last tab -------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsPage = $wb.Worksheets.Add($null, $lastSheet)
$wsPage.Name = "nbaPage"

# Note: write A2 before A1 so the shared-string table picks up "NBA"
# ahead of "Assertions", matching the authored workbook.
$wsPage.Range("A2").Value = "NBA"
$wsPage.Range("A1").Value = "Assertions"
$wsPage.Range("A1").Font.Bold = $true

# --- nbaTeams: extend with an empty, text-formatted row 32 -------------
$wsTeams = $wb.Worksheets.Item("nbaTeams")
$wsTeams.Range("A32").NumberFormat = "@"
[void]$wsTeams.Range("A32").Select()
